$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-20 and 22 have their species/observation-specific data (columns
# A, B, D, E, F, G, H, Q, R) redistributed across rows; row 21 is untouched.
# Row 20 additionally receives a brand-new Id (111866265) rather than the
# literal Id carried by the rest of the row-6-sourced data.

# Row 3
$ws.Range("A3").Value = 111866276
$ws.Range("B3").Value = 78107
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6453
$ws.Range("F3").Value = 'Vedskivlav'
$ws.Range("G3").Value = 'Hertelidea botryosa'
$ws.Range("H3").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q3").Value = 702660.5304515015
$ws.Range("R3").Value = 7299928.856484808

# Row 4
$ws.Range("A4").Value = 111866301
$ws.Range("B4").Value = 90660
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 4362
$ws.Range("F4").Value = 'Blå taggsvamp'
$ws.Range("G4").Value = 'Hydnellum caeruleum'
$ws.Range("H4").Value = '(Hornem.) P.Karst.'
$ws.Range("Q4").Value = 702522.1051459431
$ws.Range("R4").Value = 7300047.742725079

# Row 5
$ws.Range("A5").Value = 111866194
$ws.Range("B5").Value = 90682
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 2059
$ws.Range("F5").Value = 'Skrovlig taggsvamp'
$ws.Range("G5").Value = 'Hydnellum scabrosum'
$ws.Range("H5").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q5").Value = 702686.7518818546
$ws.Range("R5").Value = 7299919.985876646

# Row 6
$ws.Range("A6").Value = 111866021
$ws.Range("B6").Value = 78107
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = 'Vedskivlav'
$ws.Range("G6").Value = 'Hertelidea botryosa'
$ws.Range("H6").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q6").Value = 702738.1111920479
$ws.Range("R6").Value = 7299806.49869829

# Row 7
$ws.Range("A7").Value = 111866048
$ws.Range("B7").Value = 90682
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 2059
$ws.Range("F7").Value = 'Skrovlig taggsvamp'
$ws.Range("G7").Value = 'Hydnellum scabrosum'
$ws.Range("H7").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q7").Value = 702750.1350314748
$ws.Range("R7").Value = 7299799.924799141

# Row 8
$ws.Range("A8").Value = 111865263
$ws.Range("B8").Value = 90658
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 4361
$ws.Range("F8").Value = 'Orange taggsvamp'
$ws.Range("G8").Value = 'Hydnellum aurantiacum'
$ws.Range("H8").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q8").Value = 702714.1819675351
$ws.Range("R8").Value = 7299724.394724619

# Row 9
$ws.Range("A9").Value = 111866170
$ws.Range("B9").Value = 90682
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 2059
$ws.Range("F9").Value = 'Skrovlig taggsvamp'
$ws.Range("G9").Value = 'Hydnellum scabrosum'
$ws.Range("H9").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q9").Value = 702754.3208386695
$ws.Range("R9").Value = 7299886.818591502

# Row 10
$ws.Range("A10").Value = 111866031
$ws.Range("B10").Value = 78107
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6453
$ws.Range("F10").Value = 'Vedskivlav'
$ws.Range("G10").Value = 'Hertelidea botryosa'
$ws.Range("H10").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q10").Value = 702750.1350314748
$ws.Range("R10").Value = 7299799.924799141

# Row 11
$ws.Range("A11").Value = 111865488
$ws.Range("B11").Value = 90660
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 4362
$ws.Range("F11").Value = 'Blå taggsvamp'
$ws.Range("G11").Value = 'Hydnellum caeruleum'
$ws.Range("H11").Value = '(Hornem.) P.Karst.'
$ws.Range("Q11").Value = 702716.2360189059
$ws.Range("R11").Value = 7299724.539719297

# Row 12
$ws.Range("A12").Value = 111865961
$ws.Range("B12").Value = 77267
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6446
$ws.Range("F12").Value = 'Kolflarnlav'
$ws.Range("G12").Value = 'Carbonicola anthracophila'
$ws.Range("H12").Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range("Q12").Value = 702714.4770808229
$ws.Range("R12").Value = 7299790.39698876

# Row 13
$ws.Range("A13").Value = 111865981
$ws.Range("B13").Value = 90652
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 3100
$ws.Range("F13").Value = 'Talltaggsvamp'
$ws.Range("G13").Value = 'Bankera fuligineoalba'
$ws.Range("H13").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q13").Value = 702695.6801449896
$ws.Range("R13").Value = 7299770.100652335

# Row 14
$ws.Range("A14").Value = 111865866
$ws.Range("B14").Value = 90652
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 3100
$ws.Range("F14").Value = 'Talltaggsvamp'
$ws.Range("G14").Value = 'Bankera fuligineoalba'
$ws.Range("H14").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q14").Value = 702753.3055412351
$ws.Range("R14").Value = 7299801.798166115

# Row 15
$ws.Range("A15").Value = 111865919
$ws.Range("B15").Value = 95538
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 221941
$ws.Range("F15").Value = 'Plattlummer'
$ws.Range("G15").Value = 'Lycopodium complanatum'
$ws.Range("H15").Value = 'L.'
$ws.Range("Q15").Value = 702755.0230470664
$ws.Range("R15").Value = 7299754.083126943

# Row 16
$ws.Range("A16").Value = 111866065
$ws.Range("B16").Value = 78107
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 6453
$ws.Range("F16").Value = 'Vedskivlav'
$ws.Range("G16").Value = 'Hertelidea botryosa'
$ws.Range("H16").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q16").Value = 702767.9701038125
$ws.Range("R16").Value = 7299827.988589783

# Row 17
$ws.Range("A17").Value = 111866131
$ws.Range("B17").Value = 90682
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 2059
$ws.Range("F17").Value = 'Skrovlig taggsvamp'
$ws.Range("G17").Value = 'Hydnellum scabrosum'
$ws.Range("H17").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q17").Value = 702756.5806601554
$ws.Range("R17").Value = 7299854.813386399

# Row 18
$ws.Range("A18").Value = 111866159
$ws.Range("B18").Value = 90652
$ws.Range("D18").Value = 'NT'
$ws.Range("E18").Value = 3100
$ws.Range("F18").Value = 'Talltaggsvamp'
$ws.Range("G18").Value = 'Bankera fuligineoalba'
$ws.Range("H18").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("Q18").Value = 702755.4455659754
$ws.Range("R18").Value = 7299865.042498757

# Row 19
$ws.Range("A19").Value = 111865668
$ws.Range("B19").Value = 78107
$ws.Range("D19").Value = 'NT'
$ws.Range("E19").Value = 6453
$ws.Range("F19").Value = 'Vedskivlav'
$ws.Range("G19").Value = 'Hertelidea botryosa'
$ws.Range("H19").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q19").Value = 702740.9003275807
$ws.Range("R19").Value = 7299743.601162716

# Row 20
$ws.Range("A20").Value = 111866265
$ws.Range("B20").Value = 78107
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 6453
$ws.Range("F20").Value = 'Vedskivlav'
$ws.Range("G20").Value = 'Hertelidea botryosa'
$ws.Range("H20").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("Q20").Value = 702680.6244306123
$ws.Range("R20").Value = 7299924.914052285

# Row 22
$ws.Range("A22").Value = 111865578
$ws.Range("B22").Value = 90854
$ws.Range("D22").Value = 'NT'
$ws.Range("E22").Value = 2079
$ws.Range("F22").Value = 'Nordtagging'
$ws.Range("G22").Value = 'Odonticium romellii'
$ws.Range("H22").Value = '(S.Lundell) Parmasto'
$ws.Range("Q22").Value = 702741.9879008483
$ws.Range("R22").Value = 7299745.739876431

# Row 11 loses its public-comment note (it now carries row 13's data);
# row 17 gains the note that used to live on row 5.
$ws.Range("AC11").ClearContents()
$ws.Range("AC17").Value = 'Flera fruktkoppar som växer i en häxring'

Write-Host "Row reshuffle applied"
